# Update coin price table rows with refreshed data from the latest scrape.
# Columns: B=Coin name, C=Link, D=Price (stored as text), E=Volume(1h) label.
# A (index), F (date) and G (hour) are unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 'BNB'
$ws.Range("C2").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D2").Value = "'243.50"
$ws.Range("E2").Value = '1BNBBNB'

$ws.Range("B3").Value = 'OKB'
$ws.Range("C3").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D3").Value = "'23.20"
$ws.Range("E3").Value = '2OKBOKB'

$ws.Range("B4").Value = 'HuobiToken'
$ws.Range("C4").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D4").Value = "'5.394"
$ws.Range("E4").Value = '3HuobiTokenHT'

$ws.Range("B5").Value = 'Cronos'
$ws.Range("C5").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D5").Value = "'0.05974"
$ws.Range("E5").Value = '4CronosCRO'

$ws.Range("B6").Value = 'GateToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D6").Value = "'3.438"
$ws.Range("E6").Value = '5GateTokenGT'

$ws.Range("B7").Value = 'KuCoinToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D7").Value = "'6.529"
$ws.Range("E7").Value = '6KuCoinTokenKCS'

$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = "'0.8104"
$ws.Range("E8").Value = '7MXTokenMX'

$ws.Range("B9").Value = 'FTXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D9").Value = "'0.9272"
$ws.Range("E9").Value = '8FTXTokenFTT'

$ws.Range("B10").Value = 'One'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D10").Value = "'0.01114"
$ws.Range("E10").Value = '9OneONEBestin24h'

$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = "'0.1423"
$ws.Range("E11").Value = '10WazirXWRX'

$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = "'0.07415"
$ws.Range("E12").Value = '11MandalaExchangeTokenMDX'

$ws.Range("B13").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C13").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D13").Value = "'0.03309"
$ws.Range("E13").Value = '12LiechtensteinCryptoassetsExchangeLCX'

$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D14").Value = "'0.03080"
$ws.Range("E14").Value = '13BitrueCoinBTR'

$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").Value = "'0.09355"
$ws.Range("E15").Value = '14BitMartTokenBMX'

$ws.Range("B16").Value = 'MCDex'
$ws.Range("C16").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D16").Value = "'3.854"
$ws.Range("E16").Value = '15MCDexMCB'

$ws.Range("B17").Value = 'BitForexToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D17").Value = "'0.001574"
$ws.Range("E17").Value = '16BitForexTokenBF'

$ws.Range("B18").Value = 'CoinExToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D18").Value = "'0.04724"
$ws.Range("E18").Value = '17CoinExTokenCET'

$ws.Range("B19").Value = 'TigerCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D19").Value = "'0.005873"
$ws.Range("E19").Value = '18TigerCashTCH'

$ws.Range("B21").Value = 'HotbitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D21").Value = "'0.004892"
$ws.Range("E21").Value = '20HotbitTokenHTB'

$ws.Range("B22").Value = 'NitroEx'
$ws.Range("C22").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D22").Value = "'0.00006803"
$ws.Range("E22").Value = '21NitroExNTX'

$ws.Range("B23").Value = 'LEO'
$ws.Range("C23").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D23").Value = "'3.570"
$ws.Range("E23").Value = '22LEOLEO'

$ws.Range("B27").Value = 'UpBots'
$ws.Range("C27").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D27").Value = "'0.0002340"
$ws.Range("E27").Value = '26UpBotsUBXT'

$ws.Range("B40").Value = 'IDEX'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ZiRElvGxqQaf+idex-idex'
$ws.Range("D40").Value = "'0.03968"
$ws.Range("E40").Value = '39IDEXIDEX'

$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D41").Value = "'0.006432"
$ws.Range("E41").Value = '40KickTokenKICK'

$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D42").Value = "'0.004202"
$ws.Range("E42").Value = '41CEJICEJI'

$ws.Range("B43").Value = 'BKEXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D43").Value = "'0.1077"
$ws.Range("E43").Value = '42BKEXTokenBKK'

$ws.Range("B44").Value = 'LocalTraders'
$ws.Range("C44").Value = 'https://coinranking.com/coin/E6DwMU2zXb+localtraders-lct'
$ws.Range("D44").Value = "'0.009202"
$ws.Range("E44").Value = '43LocalTradersLCT'

$ws.Range("B45").Value = 'CoinLion'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sot4vgRyjNXek+coinlion-lion'
$ws.Range("D45").Value = "'0.00005059"
$ws.Range("E45").Value = '44CoinLionLION'

$ws.Range("B46").Value = 'Kangarootoken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/zkVNkSGwZ3+kangarootoken-gar'
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("E46").Value = '45KangarootokenGAR'

$ws.Range("B47").Value = 'CoinbaseStockToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range("D47").Value = "'0.7002"
$ws.Range("E47").Value = '46CoinbaseStockTokenCOINWorstin24h'

$ws.Range("B48").Value = 'BOLO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range("D48").Value = "'0.002412"
$ws.Range("E48").Value = '47BOLOBOLO'

$ws.Range("B49").Value = 'CryptobidCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/h39bvStAP+cryptobidcoin-cbc'
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("E49").Value = '48CryptobidCoinCBC'
